$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row for transport costs
$ws.Range("A8").Value = "TRANSPORTE"

# Fix existing "Refeiçoes" label to "REFEICOES" (row 7, col A)
$ws.Range("A7").Value = "REFEICOES"

$ws.Range("B8").Value = "R$ 49.780"

# Update selection to reflect new active cell after edits
$ws.Range("B9").Select()
